$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# --- Correct the existing week's Riders / Average figures (Pilot Target column unchanged) ---
$ws.Range("C2").Value = 229
$ws.Range("D2").Value = 104.35

$ws.Range("C3").Value = 219
$ws.Range("D3").Value = 108.89

$ws.Range("C4").Value = 222
$ws.Range("D4").Value = 114.05

$ws.Range("C5").Value = 191
$ws.Range("D5").Value = 111.76

$ws.Range("C6").Value = 260
$ws.Range("D6").Value = 109.94

# --- Append the weekend rows (Saturday 24 Sep, Sunday 25 Sep) ---
$ws.Range("A7").Value = "Saturday"
$ws.Range("B7").Value = "24 Sep 2016"
$ws.Range("C7").Value = 107
$ws.Range("D7").Value = 50.87
$ws.Range("E7").Value = 82.36

$ws.Range("A8").Value = "Sunday"
$ws.Range("B8").Value = "25 Sep 2016"
$ws.Range("C8").Value = 67
$ws.Range("D8").Value = 39.33
$ws.Range("E8").Value = 82.53

# --- Extend the chart's plotted ranges down through the two new rows ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$ser1 = $chart.SeriesCollection(1)
$ser1.Formula = '=SERIES("Ridership",Ridership!$B$2:$B$8,Ridership!$C$2:$C$8,1)'

$ser2 = $chart.SeriesCollection(2)
$ser2.Formula = '=SERIES("Average",Ridership!$B$2:$B$8,Ridership!$D$2:$D$8,2)'

$ser3 = $chart.SeriesCollection(3)
$ser3.Formula = '=SERIES("Pilot",Ridership!$B$2:$B$8,Ridership!$E$2:$E$8,3)'

# --- Move the chart down two rows to keep it below the (now longer) table ---
$co.Top = $co.Top + (2 * $ws.Rows.RowHeight)
